# Rerun with updated histologies: update p-values on several
# per-histology sheets in the demo-clin-pvalues-by-histology workbook.

$wb = $excel.ActiveWorkbook

function Set-PValues {
    param(
        [string]$SheetName,
        [hashtable]$Updates
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($cell in $Updates.Keys) {
        $ws.Range($cell).Value = $Updates[$cell]
    }
}

Set-PValues "Low-grade glioma" @{
    "B2" = 0.258370854717028
    "B3" = 0.1408182448886
    "B4" = 0.830977487111851
    "B5" = 0.661238327317944
    "B6" = 0.143818354993677
    "B7" = 0.868300106934329
    "B8" = 0.148440568057781
    "B9" = 0.232016027318013
}

Set-PValues "Medulloblastoma" @{
    "B9" = 0.0243140925220621
}

Set-PValues "Mixed neuronal-glial tumor" @{
    "B9" = 0.434462529964249
}

Set-PValues "Ependymoma" @{
    "B9" = 0.33110233659148
}

Set-PValues "Other high-grade glioma" @{
    "B9" = 0.589401718172365
}

Set-PValues "Craniopharyngioma" @{
    "B9" = 0.882365751773955
}

Set-PValues "ATRT" @{
    "B9" = 0.160726716942899
}

Set-PValues "Meningioma" @{
    "B9" = 0.521993004459526
}

Set-PValues "DIPG or DMG" @{
    "B9" = 0.824214113519518
}

Set-PValues "Mesenchymal tumor" @{
    "B3" = 0.335913312693498
    "B9" = 0.880504238997841
}

Set-PValues "Neurofibroma plexiform" @{
    "B2" = 0.43956043956044
    "B6" = 0.571428571428572
    "B7" = 0.373626373626374
    "B8" = 0.23976023976024
    "B9" = 0.394496316665903
}

Set-PValues "Non-neoplastic tumor" @{
    "B2" = 0.708636291244987
    "B3" = 0.685990338164251
    "B4" = 0.227445018749371
    "B5" = 0.522998740390055
    "B6" = 0.120821207777732
    "B7" = 0.417317976141506
    "B8" = 0.600884409707939
    "B9" = 0.530489581031042
}

Set-PValues "Other tumor" @{
    "B9" = 0.446428676155261
}
